$wb = $excel.ActiveWorkbook

# ALC row 33 (item id 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 162.11765
$ws.Range("I33").Value = 105.22222
$ws.Range("J33").Value = 226.125
$ws.Range("K33").Value = 105.22222
$ws.Range("L33").Value = 226.125
$ws.Range("M33").Value = 123.77778
$ws.Range("N33").Value = -684.125

# ALC row 107 (item id 27766)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 378.41666
$ws.Range("I107").Value = 376.45456
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 376.45456
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1543.54544
$ws.Range("N107").Value = -4240

# ALC row 113 (item id 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3250.6875
$ws.Range("I113").Value = 3327.7273
$ws.Range("J113").Value = 3081.2
$ws.Range("K113").Value = 3327.7273
$ws.Range("L113").Value = 3081.2
$ws.Range("M113").Value = -73.72730000000001
$ws.Range("N113").Value = -9589.200000000001

# ALC row 132 (item id 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3065317
$ws.Range("I132").Value = 3004.7334
$ws.Range("J132").Value = 49000000
$ws.Range("K132").Value = 9014.200199999999
$ws.Range("L132").Value = 147000000
$ws.Range("M132").Value = -6484.200199999999
$ws.Range("N132").Value = -147005060

# ARM row 32 (item id 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20301.967
$ws.Range("I32").Value = 22756.941
$ws.Range("J32").Value = 4651.5
$ws.Range("K32").Value = 22756.941
$ws.Range("L32").Value = 4651.5
$ws.Range("M32").Value = -22469.941
$ws.Range("N32").Value = -5225.5

# ARM row 61 (item id 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 111334710
$ws.Range("I61").Value = 143001730
$ws.Range("J61").Value = 500150
$ws.Range("K61").Value = 143001730
$ws.Range("L61").Value = 500150
$ws.Range("M61").Value = -143001518
$ws.Range("N61").Value = -500574

# ARM row 74 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4422831
$ws.Range("I74").Value = 5579217.5
$ws.Range("K74").Value = 5579217.5
$ws.Range("M74").Value = -5578343.5

# ARM row 77 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4422831
$ws.Range("I77").Value = 5579217.5
$ws.Range("K77").Value = 27896087.5
$ws.Range("M77").Value = -27891719.5

# ARM row 122 (item id 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 55557056
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -333338230

# ARM row 132 (item id 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 73587.17
$ws.Range("I132").Value = 46383.184
$ws.Range("K132").Value = 139149.552
$ws.Range("M132").Value = -136619.552

# ARM row 136 (item id 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 111334710
$ws.Range("I136").Value = 143001730
$ws.Range("J136").Value = 500150
$ws.Range("K136").Value = 429005190
$ws.Range("L136").Value = 1500450
$ws.Range("M136").Value = -429002640
$ws.Range("N136").Value = -1505550

# BSM row 20 (item id 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1200.6471
$ws.Range("I20").Value = 991.0909
$ws.Range("K20").Value = 991.0909
$ws.Range("M20").Value = -744.0909

# BSM row 94 (item id 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 596.9666999999999
$ws.Range("I94").Value = 609.65
$ws.Range("J94").Value = 571.6
$ws.Range("K94").Value = 609.65
$ws.Range("L94").Value = 571.6
$ws.Range("M94").Value = -158.65
$ws.Range("N94").Value = -1473.6

# BSM row 107 (item id 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2125.5715
$ws.Range("I107").Value = 2013.75
$ws.Range("J107").Value = 2274.6667
$ws.Range("K107").Value = 2013.75
$ws.Range("L107").Value = 2274.6667
$ws.Range("M107").Value = -93.75
$ws.Range("N107").Value = -6114.6667

# BSM row 134 (item id 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3902.8333
$ws.Range("I134").Value = 2886.9092
$ws.Range("J134").Value = 5499.2856
$ws.Range("K134").Value = 8660.7276
$ws.Range("L134").Value = 16497.8568
$ws.Range("M134").Value = -6125.7276
$ws.Range("N134").Value = -21567.8568

# CRP row 99 (item id 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2662.4
$ws.Range("I99").Value = 2210.5881
$ws.Range("K99").Value = 2210.5881
$ws.Range("M99").Value = -712.5880999999999

# CRP row 107 (item id 27689)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 330.80768
$ws.Range("I107").Value = 297
$ws.Range("J107").Value = 406.875
$ws.Range("K107").Value = 297
$ws.Range("L107").Value = 406.875
$ws.Range("M107").Value = 1623
$ws.Range("N107").Value = -4246.875

# CRP row 126 (item id 36198)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2662.4
$ws.Range("I126").Value = 2210.5881
$ws.Range("K126").Value = 6631.7643
$ws.Range("M126").Value = -4161.7643

# CRP row 132 (item id 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 22396.307
$ws.Range("I132").Value = 1725.6487
$ws.Range("J132").Value = 86130.836
$ws.Range("K132").Value = 5176.9461
$ws.Range("L132").Value = 258392.508
$ws.Range("M132").Value = -2646.9461
$ws.Range("N132").Value = -263452.508

# CRP row 134 (item id 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 25243.916
$ws.Range("I134").Value = 1881.1428
$ws.Range("J134").Value = 188783.33
$ws.Range("K134").Value = 5643.428400000001
$ws.Range("L134").Value = 566349.99
$ws.Range("M134").Value = -3108.428400000001
$ws.Range("N134").Value = -571419.99

# GSM row 75 (item id 11008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 45916.668
$ws.Range("J75").Value = 45916.668
$ws.Range("L75").Value = 45916.668
$ws.Range("N75").Value = -47664.668

# GSM row 78 (item id 11008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H78").Value = 45916.668
$ws.Range("J78").Value = 45916.668
$ws.Range("L78").Value = 137750.004
$ws.Range("N78").Value = -146486.004

# GSM row 113 (item id 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1168.6666
$ws.Range("I113").Value = 865.3125
$ws.Range("J113").Value = 1775.375
$ws.Range("K113").Value = 865.3125
$ws.Range("L113").Value = 1775.375
$ws.Range("M113").Value = 1304.6875
$ws.Range("N113").Value = -6115.375

# GSM row 132 (item id 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 82101.36
$ws.Range("I132").Value = 68251.734
$ws.Range("J132").Value = 102875.8
$ws.Range("K132").Value = 204755.202
$ws.Range("L132").Value = 308627.4
$ws.Range("M132").Value = -202225.202
$ws.Range("N132").Value = -313687.4

# LTW row 16 (item id 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3010.5
$ws.Range("I16").Value = 1095.5834
$ws.Range("J16").Value = 14500
$ws.Range("K16").Value = 1095.5834
$ws.Range("L16").Value = 14500
$ws.Range("M16").Value = -925.5834
$ws.Range("N16").Value = -14840

# LTW row 100 (item id 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 38468.777
$ws.Range("I100").Value = 67696.87
$ws.Range("J100").Value = 1933.6666
$ws.Range("K100").Value = 67696.87
$ws.Range("L100").Value = 1933.6666
$ws.Range("M100").Value = -67155.87
$ws.Range("N100").Value = -3015.6666

# LTW row 122 (item id 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2333.3333
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").ClearContents()

# LTW row 132 (item id 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 58637.777
$ws.Range("I132").Value = 3016.5454
$ws.Range("J132").Value = 146042.58
$ws.Range("K132").Value = 9049.636200000001
$ws.Range("L132").Value = 438127.74
$ws.Range("M132").Value = -6519.636200000001
$ws.Range("N132").Value = -443187.74

# WVR row 122 (item id 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2246.9524
$ws.Range("I122").Value = 969.38464
$ws.Range("K122").Value = 2908.15392
$ws.Range("M122").Value = -458.1539199999997

# WVR row 126 (item id 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2295
$ws.Range("I126").Value = 1774.909
$ws.Range("K126").Value = 5324.727000000001
$ws.Range("M126").Value = -2854.727000000001

# WVR row 132 (item id 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 47598.906
$ws.Range("I132").Value = 30178.234
$ws.Range("J132").Value = 113410.336
$ws.Range("K132").Value = 90534.702
$ws.Range("L132").Value = 340231.008
$ws.Range("M132").Value = -88004.702
$ws.Range("N132").Value = -345291.008

# WVR row 136 (item id 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 41271.12
$ws.Range("I136").Value = 23123.467
$ws.Range("K136").Value = 69370.401
$ws.Range("M136").Value = -66820.401
